$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "Hello"
$ws.Range("A1").Font.Name = "Calibri"
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Bold = $false
$ws.Range("A1").Font.FontStyle = "Regular"
